# live_trading_results.xlsx update
# Trade #71 (MarketMaking, opened 2026-02-17 20:48:50) closes at 21:04:20
# (early_exit), and a brand-new trade #104 is opened at 21:04:08.
# Summary / Strategy Status roll-up metrics are refreshed accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1400.58   # Current Capital
$summary.Range("B4").Value = 0.38      # Total P&L $
$summary.Range("B5").Value = 0.11      # Total P&L %
$summary.Range("B6").Value = 71        # Total Trades
$summary.Range("B7").Value = 33        # Winning Trades
$summary.Range("B9").Value = 46.48     # Win Rate %

# ---------------------------------------------------------------------
# 2) Strategy Status sheet - MarketMaking row
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 100.58     # Capital
$status.Range("D5").Value = 38         # Trades
$status.Range("E5").Value = 0.27       # P&L $
$status.Range("F5").Value = 0.58       # P&L %
$status.Range("G5").Value = 50         # Win Rate %

# ---------------------------------------------------------------------
# 3) All Trades sheet
#    - row 72 = Trade #71 -> close it out
#    - row 105 = new Trade #104 (freshly opened)
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

$allTrades.Cells.Item(72, 7).Value = 0.99          # G72 Exit Price
$allTrades.Cells.Item(72, 8).Value = "CLOSED"      # H72 Status
$allTrades.Cells.Item(72, 9).Value = 2.0619        # I72 P&L %
$allTrades.Cells.Item(72, 10).Value = 0.02         # J72 P&L $
$allTrades.Cells.Item(72, 11).Value = 100.58       # K72 Capital After
$allTrades.Cells.Item(72, 12).Value = "early_exit" # L72 Exit Reason
$allTrades.Cells.Item(72, 13).Value = 0.23         # M72 Duration (min)

$allTrades.Cells.Item(105, 1).Value = 104                                    # A105 Trade #
$allTrades.Cells.Item(105, 2).NumberFormat = "@"
$allTrades.Cells.Item(105, 2).Value = "2026-02-17"                           # B105 Date (force text, else Excel parses as a date serial)
$allTrades.Cells.Item(105, 2).ClearFormats()                                 # drop the temporary text format again
$allTrades.Cells.Item(105, 3).Value = "21:04:08"                             # C105 Time
$allTrades.Cells.Item(105, 4).Value = "MarketMaking"                         # D105 Strategy
$allTrades.Cells.Item(105, 5).Value = "DOWN"                                 # E105 Side
$allTrades.Cells.Item(105, 6).Value = 0.97                                   # F105 Entry Price
$allTrades.Cells.Item(105, 8).Value = "OPEN"                                 # H105 Status
$allTrades.Cells.Item(105, 9).Value = 0                                      # I105 P&L %
$allTrades.Cells.Item(105, 10).Value = 0                                     # J105 P&L $
$allTrades.Cells.Item(105, 11).Value = 100.5619219857093                     # K105 Capital After
$allTrades.Cells.Item(105, 13).Value = 0                                     # M105 Duration (min)
$allTrades.Cells.Item(105, 14).Value = 0                                     # N105 Entry Slippage (bps)
$allTrades.Cells.Item(105, 15).Value = 0                                     # O105 Exit Slippage (bps)
$allTrades.Cells.Item(105, 16).Value = 0.6                                   # P105 Confidence
$allTrades.Cells.Item(105, 17).Value = "Normal spread capture: 19600 bps"    # Q105 Entry Reason

# ---------------------------------------------------------------------
# 4) MarketMaking sheet
#    - row 39 = Trade #71 -> close it out
#    - row 72 = new Trade #104 (freshly opened)
# ---------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")

$mm.Cells.Item(39, 7).Value = 0.99          # G39 Exit Price
$mm.Cells.Item(39, 8).Value = "CLOSED"      # H39 Status
$mm.Cells.Item(39, 9).Value = 2.0619        # I39 P&L %
$mm.Cells.Item(39, 10).Value = 0.02         # J39 P&L $
$mm.Cells.Item(39, 11).Value = 100.58       # K39 Capital After
$mm.Cells.Item(39, 16).Value = "early_exit" # P39 Exit Reason
$mm.Cells.Item(39, 17).Value = 0.23         # Q39 Duration (min)

$mm.Cells.Item(72, 1).Value = 104                                  # A72 Trade #
$mm.Cells.Item(72, 2).NumberFormat = "@"
$mm.Cells.Item(72, 2).Value = "2026-02-17"                         # B72 Date (force text, else Excel parses as a date serial)
$mm.Cells.Item(72, 2).ClearFormats()                                # drop the temporary text format again
$mm.Cells.Item(72, 3).Value = "21:04:08"                           # C72 Time
$mm.Cells.Item(72, 4).Value = "MarketMaking"                       # D72 Strategy
$mm.Cells.Item(72, 5).Value = "DOWN"                                # E72 Side
$mm.Cells.Item(72, 6).Value = 0.97                                  # F72 Entry Price
$mm.Cells.Item(72, 8).Value = "OPEN"                                # H72 Status
$mm.Cells.Item(72, 9).Value = 0                                     # I72 P&L %
$mm.Cells.Item(72, 10).Value = 0                                    # J72 P&L $
$mm.Cells.Item(72, 11).Value = 100.5619219857093                    # K72 Capital After
$mm.Cells.Item(72, 12).Value = 0                                    # L72 Entry Slippage (bps)
$mm.Cells.Item(72, 13).Value = 0                                    # M72 Exit Slippage (bps)
$mm.Cells.Item(72, 14).Value = 0.6                                  # N72 Confidence
$mm.Cells.Item(72, 15).Value = "Normal spread capture: 19600 bps"   # O72 Entry Reason
$mm.Cells.Item(72, 17).Value = 0                                    # Q72 Duration (min)
